$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from its old spot (just before the
#    "Jan/25 - Feb /25" run) to right after the "EXPERIENCE" heading.
# ------------------------------------------------------------------

# Remove the existing _GoBack bookmark (near "Jan/25 - Feb /25").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Find the end of the "EXPERIENCE" run.
$expRng = $d.Content
$expRng.Find.Execute("EXPERIENCE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$expPos = $expRng.End

# A zero-length Range passed to Bookmarks.Add is mishandled by this
# engine (it resets to the start of the document), so insert a throw-
# away marker character, wrap the bookmark around that single
# character, then delete the character back out again -- the
# bookmark's start/end tags stay put, collapsed at the right spot.
$insertPt = $d.Range($expPos, $expPos)
$insertPt.InsertAfter("X")
$markRng = $d.Range($expPos, $expPos + 1)
$d.Bookmarks.Add("_GoBack", $markRng)
$d.Range($expPos, $expPos + 1).Text = ""

# ------------------------------------------------------------------
# 2) Merge the "Key Projects" + ":" runs into a single run reading
#    "Key Projects:", bump the paragraph-mark size to 18pt (sz 36)
#    and the run's own size to 12pt (sz 24).
# ------------------------------------------------------------------

# Re-typing the same text over the found range merges the two runs
# that used to hold "Key Projects" and ":" into one run.
$d.Content.Find.Execute("Key Projects:", $true, $false, $false, $false, $false, $true, 1, $false, "Key Projects:", 2)

# Set the whole paragraph (text + paragraph mark) to 18pt first, this
# stamps the paragraph-mark run properties (w:pPr/w:rPr) with sz 36.
$kpRng = $d.Content
$kpRng.Find.Execute("Key Projects:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$kpPara = $kpRng.Paragraphs(1)
$kpPara.Range.Font.Size = 18

# Now shrink just the run text back down to 12pt (sz 24), leaving the
# paragraph mark's own size at 18pt (sz 36).
$kpTextRng = $d.Content
$kpTextRng.Find.Execute("Key Projects:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$kpTextRng.Font.Size = 12

Write-Host "done"
